# Applies the edit described by the diff:
#  1. Adds a first-line indent to the "Пам'ятка вкрита нестійкими ..." paragraph
#     (section 4), and
#  2. Inserts a large new "page two" section (sections 5 and 6 of the
#     restoration passport) right after the table that ends with
#     "Керівник: ..." and before the pre-existing page-break paragraph.

$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "вкрита нестійкими") {
        $p.Range.ParagraphFormat.FirstLineIndent = 18
        break
    }
}

# --- Change 2 --------------------------------------------------------------
# Locate the paragraph that holds the page break immediately following the
# table whose last cell reads "Керівник: ...".
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "Керівник") {
        $anchorIndex = $i
        break
    }
}

$pageBreakIndex = 0
for ($i = $anchorIndex; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.WordOpenXML -match 'w:type="page"') {
        $pageBreakIndex = $i
        break
    }
}

# Make room: insert one empty paragraph right before that page break, then
# replace its content with the whole new block in one shot via InsertXML.
# (The four headings that use the "Emphasis" character style are inserted as
# plain runs here and get the rStyle applied in a second pass below, because
# InsertXML does not reliably keep rStyle references on brand-new runs.)
$targetP = $d.Paragraphs($pageBreakIndex)
$targetP.Range.InsertParagraphBefore()
$newP = $d.Paragraphs($pageBreakIndex)

$xmlPayload = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>5. Основні дані з історії пам'ятки (довідка про побутування; відомості про умови зберігання, попередні дослідження, консерваційно-реставраційні заходи тощо), джерело надходження інформації</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:jc w:val="left"/></w:pPr><w:r><w:t>Пам'ятка походить з (вказати регіон) та перебувала у фондах (введіть назву установи або власність колекції).</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>6. Стан пам'ятки до реставрації</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:t>6.1 За візуальним спостереженням:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:jc w:val="left"/></w:pPr><w:r><w:t>За візуальним спостереженням: ...</w:t><w:br/><w:br/><w:t>I. Візуальне дослідження (опис пам’ятки):</w:t><w:br/><w:br/><w:t xml:space="preserve">    1.Вказати назву (якщо є спеціальний термін).</w:t><w:br/><w:br/><w:t xml:space="preserve">    2.Описати форму та колір.</w:t><w:br/><w:br/><w:t>II. Описати вигляд предмета:</w:t><w:br/><w:br/><w:t xml:space="preserve">    1. Складові предмета, їх геометрична форма;</w:t><w:br/><w:br/><w:t>III. Забруднення:</w:t><w:br/><w:br/><w:t xml:space="preserve">    1.Нестійкі (пилові, брудові, ґрунтові) .</w:t><w:br/><w:br/><w:t xml:space="preserve">    2.Стійкі (вапнякові, природні та синтетичні смоли, висоли, гіпсові забруднення, плями кислів металів, сліди кіптяви,</w:t><w:br/><w:br/><w:t xml:space="preserve">        пеку, жиру, плями від пластиліну, масляної фарби, воску, клейові забруднення, чорнила, туш,</w:t><w:br/><w:br/><w:t xml:space="preserve">        записи фарбами (якого кольору), забруднення фарбами  від попередніх тонувань – місцезнаходження фарби,</w:t><w:br/><w:br/><w:t xml:space="preserve">        забруднення на зламах фрагментів (від клею, вапнякових нашарувань, пило брудові, ґрунтові і т.д.).</w:t><w:br/><w:br/><w:t xml:space="preserve">    3.Визначити за візуальним спостереженням яким клеєм склеєно фрагменти (клеєм БФ (світло-коричневого, коричневого, червоного кольору, прозорий),</w:t><w:br/><w:br/><w:t xml:space="preserve">        ПВА (полівінилацетатний клей молочного кольору, непрозорий, безбарвний, прозорий)</w:t><w:br/><w:br/><w:t>Визначити форму забруднення (у вигляді локальних плям, неправильної форми, повсюдно, забруднення якоїсь частини пам’ятки).</w:t><w:br/><w:br/><w:t xml:space="preserve">    Матеріали:</w:t><w:br/><w:br/><w:t>(ДЕРЕВИНА):</w:t><w:br/><w:br/><w:t xml:space="preserve">    1.Вказати спосіб виготовлення (різьблення, слюсарні роботи); </w:t><w:br/><w:br/><w:t xml:space="preserve">    2.Вказати якій породі деревини відповідає орнамент.</w:t><w:br/><w:br/><w:t xml:space="preserve">    3.Описати різьблення, орнамент (заглиблений, рельєфний).</w:t><w:br/><w:br/><w:t xml:space="preserve">    4.Описати стан фарбового шару (ступінь зчеплення), розпис (монохромний, поліхромний, колір фарб)</w:t><w:br/><w:br/><w:t>(МЕТАЛ):</w:t><w:br/><w:br/><w:t xml:space="preserve">    1.Вказати спосіб виготовлення (лиття, слюсарні роботи, гравіювання, травлення, інкрустація, зернь, паяння)</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:t>6.1.2 Втрати та пошкодження:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:jc w:val="left"/></w:pPr><w:r><w:t>ІV. Попередня реставрація:</w:t><w:br/><w:br/><w:t xml:space="preserve">    1. Якщо була – описати якість попередньої реставрації).</w:t><w:br/><w:br/><w:t xml:space="preserve">    2. Свідчення про попередню реставрацію (відсутні, якщо є вказати джерело чи з чиїх слів записано).</w:t><w:br/><w:br/><w:t xml:space="preserve">    3. Реставрація не повна якщо:</w:t><w:br/><w:br/><w:t xml:space="preserve">        -  фрагменти склеєні, а втрати не восповнені;</w:t><w:br/><w:br/><w:t xml:space="preserve">        -восповнені частково;</w:t><w:br/><w:br/><w:t xml:space="preserve">        -є втрати у будь-якій частині виробу (вказати місце втрати, форму, розмір).</w:t><w:br/><w:br/><w:t>V. Опис наявних втрат та пошкоджень:</w:t><w:br/><w:br/><w:t xml:space="preserve">    1. Вказати із якої кількості фрагментів складається пам’ятка.</w:t><w:br/><w:br/><w:t xml:space="preserve">    2. Вказати на якість попереднього склеювання, доповнення.</w:t><w:br/><w:br/><w:t xml:space="preserve">    3. Вказати які частини пам’ятки відсутні (вказати розмір в см/мм та кв. см/мм).</w:t><w:br/><w:br/><w:t xml:space="preserve">    4. Вказати пошкодження, відслоюваня, розшарування, деформації.</w:t><w:br/><w:br/><w:t xml:space="preserve">    5. Визначити дефекти:</w:t><w:br/><w:br/><w:t xml:space="preserve">        -привнесені від археологічного чи реставраційного інструменту, помітки олівцем і т.д.;</w:t><w:br/><w:br/><w:t xml:space="preserve">        -виробничі;</w:t><w:br/><w:br/><w:t xml:space="preserve">        -тріщини (наскрізні, не наскрізні, волосяні (вказати форму, розмір, розташування);</w:t><w:br/><w:br/><w:t xml:space="preserve">        -сколи, вибоїни, незначні втрати, каверни, пробоїни, викришування, потертості, подряпини (вказати форму, розмір, розташування);</w:t><w:br/><w:br/><w:t xml:space="preserve">        -визначити дефекти тонувань (якщо є розпис, консерваційного покриття – описати колір, стан збереження,</w:t><w:br/><w:br/><w:t xml:space="preserve">        наявні значні або незначні, часткові, локальні втрати, потертості відшарування (вказати форму, розмір, розташування);</w:t><w:br/><w:br/><w:t>VІ. Біологічні пошкодження і руйнування:</w:t><w:br/><w:br/><w:t xml:space="preserve">    6.Бактерії, гриби та продукти їх життєдіяльності.</w:t><w:br/><w:br/><w:t>(ДЕРЕВИНА):</w:t><w:br/><w:br/><w:t xml:space="preserve">    1.Описати стан деревини (відмінний, уражений шкідниками, трухлявий, деформації, поверхня спучена, слоїста і т.д.)</w:t><w:br/><w:br/><w:t>(МЕТАЛ):</w:t><w:br/><w:br/><w:t xml:space="preserve">    1. Визначити ступінь збереженості предмета:</w:t><w:br/><w:br/><w:t xml:space="preserve">        -Предмет гарної збереженості (новий метал, локальна корозія, загальне забруднення);</w:t><w:br/><w:br/><w:t xml:space="preserve">        -Предмет задовільної збереженості (новий метал, загальне забруднення, суцільна корозія або благородна патина, </w:t><w:br/><w:br/><w:t xml:space="preserve">            декор добре читається, металеве ядро не мінералізоване,наявні потертості, незначна деформація)</w:t><w:br/><w:br/><w:t xml:space="preserve">        -Новий метал поганої збереженості (загальні пило-брудові нашарування, суцільна корозія, </w:t><w:br/><w:br/><w:t xml:space="preserve">            наявні локальні рецедивуючі продукти корозії, значна деформація, втрати).</w:t><w:br/><w:br/><w:t xml:space="preserve">        -Археологічний метал (наявний товстий шар корозійного нашарування, міжкристалітна корозія, крихкість, тріщини, </w:t><w:br/><w:br/><w:t xml:space="preserve">            втрати, крізна корозія. Декор читається погано, форми предмета та металеве ядро збережені)</w:t><w:br/><w:br/><w:t xml:space="preserve">        -Археологічний метал з частково мінералізованим металевим ядром (форма предмету читається погано, </w:t><w:br/><w:br/><w:t xml:space="preserve">            механічна міцність слабка, можливі руйнування, крізна корозія, наявна часткова мінералізація металевого ядра).</w:t><w:br/><w:br/><w:t xml:space="preserve">        -Археологічний метал з повністю мінералізованим металевим ядром (предмет перетворився на безформну масу, </w:t><w:br/><w:br/><w:t xml:space="preserve">            що складається з продуктів корозії і не має механічної міцності).</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:t>6.1.3 Старі номери та позначення:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:jc w:val="left"/></w:pPr><w:r><w:t>1. Вказати де, яким чорнилом, які номери написано (написи та шифри необхідно фотофіксувати).</w:t><w:br/><w:br/><w:t xml:space="preserve"> Марки (описати місцезнаходження, зовнішній вигляд марки або етикетки, колір і розмір, чим написано і чим приклеєний шифр).</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:t>6.1.3 Старі номери та позначення:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/><w:jc w:val="left"/></w:pPr><w:r><w:t>Довжина: 32 мм;</w:t><w:br/><w:br/><w:t>Висота: 43 мм;</w:t><w:br/><w:br/><w:t>Ширина: 5 мм;</w:t><w:br/><w:br/><w:t>Товщина: 11 мм.</w:t><w:br/><w:br/><w:t xml:space="preserve">    </w:t></w:r></w:p>
'@

[void]$newP.Range.InsertXML($xmlPayload)

# --- Fix up the four "Emphasis"-styled headings -----------------------------
$emphasisTexts = @(
    "6.1 За візуальним спостереженням:",
    "6.1.2 Втрати та пошкодження:",
    "6.1.3 Старі номери та позначення:",
    "6.1.3 Старі номери та позначення:"
)

$remaining = New-Object System.Collections.ArrayList
foreach ($t in $emphasisTexts) { [void]$remaining.Add($t) }

for ($i = $pageBreakIndex; $i -le $d.Paragraphs.Count; $i++) {
    if ($remaining.Count -eq 0) { break }
    $p = $d.Paragraphs($i)
    $ptext = $p.Range.Text
    $want = $remaining[0]
    if ($ptext -eq ($want + "`r")) {
        $scoped = $d.Range($p.Range.Start, $p.Range.End - 1)
        $scoped.Style = "Emphasis"
        $remaining.RemoveAt(0)
    }
}
